$wb = $excel.ActiveWorkbook

# --- Update Moorings sheet ---
$wsMoorings = $wb.Worksheets.Item("Moorings")
$wsMoorings.Range("A2").Value = "GI05MOAS-GL478"

# Update the selected cell / view state on the Moorings sheet
$wsMoorings.Activate()
$wsMoorings.Range("G34").Select()

# --- Update Asset_Cal_Info sheet ---
$wsAsset = $wb.Worksheets.Item("Asset_Cal_Info")
$wsAsset.Range("A2").Value = "GI05MOAS-GL478-01-FLORDM000"
$wsAsset.Range("A3").Value = "GI05MOAS-GL478-01-FLORDM000"
$wsAsset.Range("A4").Value = "GI05MOAS-GL478-01-FLORDM000"
$wsAsset.Range("A5").Value = "GI05MOAS-GL478-01-FLORDM000"
$wsAsset.Range("A7").Value = "GI05MOAS-GL478-02-DOSTAM000"
$wsAsset.Range("A9").Value = "GI05MOAS-GL478-04-CTDGVM000"
$wsAsset.Range("A11").Value = "GI05MOAS-GL478-00-ENG000000"
